$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "day" sheet: append 13 new rows (632-644) after the existing 631 rows.
# ---------------------------------------------------------------------------
$dayWs = $wb.Worksheets.Item("day")

$dayRows = @(
    @(632,  1, "OFSS",       "Oracle Financial Services Software Limited",    532466, -0.14,               11442.65, 54125,   "day", "01/10/2024 11:36:53"),
    @(633,  2, "MPHASIS",    "Mphasis Limited",                               526299,  0.17,                 3015.6, 899983,  "day", "01/10/2024 11:36:53"),
    @(634,  3, "ICICIGI",    "ICICI Lombard General Insurance Company Ltd",   540716, -0.9399999999999999,   2153.9, 381613,  "day", "01/10/2024 11:36:53"),
    @(635,  4, "SBILIFE",    "SBI Life Insurance Company Ltd",                540719, -0.53,                1834.2, 885533,  "day", "01/10/2024 11:36:53"),
    @(636,  5, "CIPLA",      "Cipla Limited",                                 500087,  0.65,                1664.85, 945740,  "day", "01/10/2024 11:36:53"),
    @(637,  6, "TECHM",      "Tech Mahindra Limited",                         532755,  3.06,                1625.4, 5176646, "day", "01/10/2024 11:36:53"),
    @(638,  7, "TATACONSUM", "TATA Consumer Products Ltd",                    500800, -0.06,                1196.25, 892383,  "day", "01/10/2024 11:36:53"),
    @(639,  8, "JINDALSTEL", "Jindal Steel & Power Limited",                  532286, -0.41,                1035.35, 2382544, "day", "01/10/2024 11:36:53"),
    @(640,  9, "HDFCLIFE",   "HDFC Life Insurance Company Ltd",               540777, -1.09,                 710.2, 1987168, "day", "01/10/2024 11:36:53"),
    @(641, 10, "BERGEPAINT", "Berger Paints (i) Limited",                     509480, -0.43,                619.65, 1098938, "day", "01/10/2024 11:36:53"),
    @(642, 11, "DABUR",      "Dabur India Limited",                           500096, -0.99,                   619, 1096213, "day", "01/10/2024 11:36:53"),
    @(643, 12, "BIOCON",     "Biocon Limited",                                532523,  2.08,                370.65, 2195461, "day", "01/10/2024 11:36:53"),
    @(644, 13, "GMRINFRA",   "Gmr Infrastructure Limited",                    532754, -0.14,   93.93000000000001, 6427524, "day", "01/10/2024 11:36:53")
)

foreach ($r in $dayRows) {
    $rowNum = $r[0]
    $dayWs.Cells.Item($rowNum, 1).Value = $r[1]
    $dayWs.Cells.Item($rowNum, 2).Value = $r[2]
    $dayWs.Cells.Item($rowNum, 3).Value = $r[3]
    $dayWs.Cells.Item($rowNum, 4).Value = $r[4]
    $dayWs.Cells.Item($rowNum, 5).Value = $r[5]
    $dayWs.Cells.Item($rowNum, 6).Value = $r[6]
    $dayWs.Cells.Item($rowNum, 7).Value = $r[7]
    $dayWs.Cells.Item($rowNum, 8).Value = $r[8]
    $dayWs.Cells.Item($rowNum, 9).Value = $r[9]
}

# ---------------------------------------------------------------------------
# 2) "month" sheet: D48:D52 were stored as text; convert to real numbers.
# ---------------------------------------------------------------------------
$monthWs = $wb.Worksheets.Item("month")

$monthCodes = @{
    48 = 500696
    49 = 532868
    50 = 500096
    51 = 524208
    52 = 532461
}

foreach ($rowNum in $monthCodes.Keys) {
    $monthWs.Cells.Item($rowNum, 4).Value = $monthCodes[$rowNum]
}
